$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the roster rows that changed (players/positions/teams were refreshed).
$ws.Range("A8").Value = "Royce O'Neale"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Phoenix Suns"

$ws.Range("A9").Value = "Keegan Murray"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Sacramento Kings"

$ws.Range("A10").Value = "Jarrett Allen"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Cleveland Cavaliers"

$ws.Range("A11").Value = "Mark Williams"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Charlotte Hornets"

$ws.Range("A12").Value = "Trey Murphy III"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "New Orleans Pelicans"

$ws.Range("A13").Value = "Jalen Duren"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Detroit Pistons"

$ws.Range("A14").Value = "Daniel Gafford"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Dallas Mavericks"

$ws.Range("A15").Value = "Karl-Anthony Towns"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "New York Knicks"

$ws.Range("A16").Value = "Jalen Johnson"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Atlanta Hawks"

# Row 19 (Jalen Duren / C / Detroit Pistons) was removed entirely - the table now ends at row 18.
$ws.Range("A19:C19").Value = $null
